$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for D, E, G columns (rows 2-51) so that numeric-looking
# strings (prices, percentages, hour numbers) are preserved as text, matching
# the original inlineStr cell type, rather than being converted to numbers.
$numRng = $ws.Range("D2:E51")
$numRng.NumberFormat = "@"
$hourRng = $ws.Range("G2:G51")
$hourRng.NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '329.00'
$ws.Range('E2').Value = '0.21%'
$ws.Range('G2').Value = '17'

# Row 3
$ws.Range('D3').Value = '44.07'
$ws.Range('E3').Value = '-0.73%'
$ws.Range('G3').Value = '17'

# Row 4
$ws.Range('D4').Value = '5.490'
$ws.Range('E4').Value = '-1.22%'
$ws.Range('G4').Value = '17'

# Row 5
$ws.Range('D5').Value = '0.08075'
$ws.Range('E5').Value = '0.08%'
$ws.Range('G5').Value = '17'

# Row 6
$ws.Range('D6').Value = '2.045'
$ws.Range('E6').Value = '6.83%'
$ws.Range('G6').Value = '17'

# Row 7
$ws.Range('B7').Value = 'MXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D7').Value = '0.9529'
$ws.Range('E7').Value = '0.01%'
$ws.Range('G7').Value = '17'

# Row 8
$ws.Range('B8').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C8').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D8').Value = '0.1121'
$ws.Range('E8').Value = '-5.80%'
$ws.Range('G8').Value = '17'

# Row 9
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').Value = '0.1887'
$ws.Range('E9').Value = '1.56%'
$ws.Range('G9').Value = '17'

# Row 10
$ws.Range('B10').Value = 'MCDex'
$ws.Range('C10').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D10').Value = '10.15'
$ws.Range('E10').Value = '1.50%'
$ws.Range('G10').Value = '17'

# Row 11
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = '0.1011'
$ws.Range('E11').Value = '3.76%'
$ws.Range('G11').Value = '17'

# Row 12
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = '0.04780'
$ws.Range('E12').Value = '9.32%'
$ws.Range('G12').Value = '17'

# Row 13
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Value = '0.1058'
$ws.Range('E13').Value = '-0.84%'
$ws.Range('G13').Value = '17'

# Row 14
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').Value = '0.001274'
$ws.Range('E14').Value = '-0.86%'
$ws.Range('G14').Value = '17'

# Row 15
$ws.Range('B15').Value = 'CoinExToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D15').Value = '0.04095'
$ws.Range('E15').Value = '-2.77%'
$ws.Range('G15').Value = '17'

# Row 16
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '0.006091'
$ws.Range('E16').Value = '1.69%'
$ws.Range('G16').Value = '17'

# Row 17
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.372'
$ws.Range('E17').Value = '-0.71%'
$ws.Range('G17').Value = '17'

# Row 18
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').Value = '4.417'
$ws.Range('E18').Value = '2.97%'
$ws.Range('G18').Value = '17'

# Row 19
$ws.Range('D19').Value = '2.621'
$ws.Range('E19').Value = '3.43%'
$ws.Range('G19').Value = '17'

# Row 20
$ws.Range('D20').Value = '0.3296'
$ws.Range('E20').Value = '-4.84%'
$ws.Range('G20').Value = '17'

# Row 21
$ws.Range('D21').Value = '0.1399'
$ws.Range('E21').Value = '-1.17%'
$ws.Range('G21').Value = '17'

# Row 22
$ws.Range('D22').Value = '0.2571'
$ws.Range('E22').Value = '2.64%'
$ws.Range('G22').Value = '17'

# Row 23
$ws.Range('D23').Value = '0.001307'
$ws.Range('E23').Value = '5.20%'
$ws.Range('G23').Value = '17'

# Row 24
$ws.Range('D24').Value = '0.004354'
$ws.Range('E24').Value = '-0.09%'
$ws.Range('G24').Value = '17'

# Row 25
$ws.Range('D25').Value = '0.0001251'
$ws.Range('E25').Value = '5.11%'
$ws.Range('G25').Value = '17'

# Row 26
$ws.Range('G26').Value = '17'

# Row 27
$ws.Range('G27').Value = '17'

# Row 28
$ws.Range('G28').Value = '17'

# Row 29
$ws.Range('G29').Value = '17'

# Row 30
$ws.Range('G30').Value = '17'

# Row 31
$ws.Range('G31').Value = '17'

# Row 32
$ws.Range('G32').Value = '17'

# Row 33
$ws.Range('G33').Value = '17'

# Row 34
$ws.Range('G34').Value = '17'

# Row 35
$ws.Range('G35').Value = '17'

# Row 36
$ws.Range('G36').Value = '17'

# Row 37
$ws.Range('G37').Value = '17'

# Row 38
$ws.Range('D38').Value = '0.02608'
$ws.Range('E38').Value = '-2.88%'
$ws.Range('G38').Value = '17'

# Row 39
$ws.Range('D39').Value = '0.05637'
$ws.Range('E39').Value = '1.75%'
$ws.Range('G39').Value = '17'

# Row 40
$ws.Range('D40').Value = '0.007573'
$ws.Range('E40').Value = '0.36%'
$ws.Range('G40').Value = '17'

# Row 41
$ws.Range('D41').Value = '0.1402'
$ws.Range('E41').Value = '-0.34%'
$ws.Range('G41').Value = '17'

# Row 42
$ws.Range('D42').Value = '0.007370'
$ws.Range('E42').Value = '-11.80%'
$ws.Range('G42').Value = '17'

# Row 43
$ws.Range('D43').Value = '0.001981'
$ws.Range('E43').Value = '-1.70%'
$ws.Range('G43').Value = '17'

# Row 44
$ws.Range('D44').Value = '0.008823'
$ws.Range('E44').Value = '-0.72%'
$ws.Range('G44').Value = '17'

# Row 45
$ws.Range('D45').Value = '0.00007075'
$ws.Range('G45').Value = '17'

# Row 46
$ws.Range('E46').Value = '0.02%'
$ws.Range('G46').Value = '17'

# Row 47
$ws.Range('E47').Value = '-0.19%'
$ws.Range('G47').Value = '17'

# Row 48
$ws.Range('E48').Value = '54.12%'
$ws.Range('G48').Value = '17'

# Row 49
$ws.Range('E49').Value = '23.73%'
$ws.Range('G49').Value = '17'

# Row 50
$ws.Range('E50').Value = '0.02%'
$ws.Range('G50').Value = '17'

# Row 51
$ws.Range('E51').Value = '0.02%'
$ws.Range('G51').Value = '17'

# Reset style so no extra explicit style index is attached to the cells
# (keeps them on the default/no style, matching the original workbook).
$numRng.Style = "Normal"
$hourRng.Style = "Normal"
